# Updated cryptos list (price/volume refresh + Aptos/BabyDogeCoin rank swap
# at rows 48-49), matching the Wed Aug 9 17:09:41 UTC 2023 GitHub Actions run.
#
# Price cells (column D) that would otherwise be auto-parsed as numbers by
# Excel are briefly forced to Text via NumberFormat="@" so they keep the
# exact original string layout (e.g. "243.37" instead of numeric 243.37),
# then the cell style is reset back to "Normal" so no stray number format
# sticks around on the cell itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.562.20'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.850.90'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6529'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.68%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.95'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07484'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2973'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.48'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07639'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '1.842.42'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.039'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6856'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009548'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.121'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").Value = '29.573.08'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '2.109.36'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '236.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.698'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.003'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1424'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.518'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06038'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.492'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.260'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.137'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.075'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.869'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7275'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.798'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").Value = '1.200.79'
$ws.Range("E41").Value = '  -2.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.275'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9102'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '2.014.90'
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.383'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.05%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000122'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4052'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.118'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.17%  '
